$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Unique" column (L) mirroring the existing "Multivalued" column (K).
# Copy K4 (bold header style) into L4, then change its text to "Unique".
$ws.Range("K4").Copy($ws.Range("L4")) | Out-Null
$ws.Range("L4").Value = "Unique"

# Copy K5:K8 (text "FALSE" cells with the TRUE/FALSE display format) into L5:L8
# unchanged, since the new "Unique" column has the same "FALSE" value for every row.
$ws.Range("K5").Copy($ws.Range("L5")) | Out-Null
$ws.Range("K6").Copy($ws.Range("L6")) | Out-Null
$ws.Range("K7").Copy($ws.Range("L7")) | Out-Null
$ws.Range("K8").Copy($ws.Range("L8")) | Out-Null

# Update the active selection to match the edited region (L7:L8).
$ws.Range("L7:L8").Select() | Out-Null
